$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 4
$ws.Range("C3").Value = 4
$ws.Range("C4").Value = 6
$ws.Range("C5").Value = 6
$ws.Range("C6").Value = 16
$ws.Range("C7").Value = 16
$ws.Range("C9").Value = 8

$ws.Range("C6").Select()

$wb.Save()
